$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: table cell "Mr Subroto Ghosh" -> two runs "Mr." + " Subroto Ghosh"
# ---------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$r1 = $d.Range($p12.Range.Start, $p12.Range.Start + 3)   # "Mr "
$r1.Text = "Mr."

$d = $word.ActiveDocument
$p12 = $d.Paragraphs.Item(12)
$r2Start = $p12.Range.Start + 3                          # right after "Mr."
$r2 = $d.Range($r2Start, $p12.Range.End)
$r2.InsertBefore(" ")

# Force a run split by toggling a formatting property on/off so that the
# two pieces of text end up as two separate (but identically formatted) runs.
$d = $word.ActiveDocument
$p12 = $d.Paragraphs.Item(12)
$r2 = $d.Range($r2Start, $p12.Range.End)
$r2.Bold = 1
$r2.Bold = 0

# ---------------------------------------------------------------------
# Edit 2: re-work the 4 "mock test process" bullet points.
#
# Before:
#   34: We will give names of students and phone number to College Doors
#   35: College Doors will register students' names in database
#   36: Faculties/CIO team will upload question paper by 24th November 2023
#   37: After completion of correction marks will be sent to students in WhatsApp
#
# After:
#   34: Faculties/CIO team will upload question paper by 24th November 2023
#   35: College Doors will create login credential for each faculties.
#   36: We will give names of students and phone number to College Doors
#   37 (new): College Doors will register students' names in database by 24th November 2023
#   38: After completion of correction marks will be sent to students in WhatsApp
# ---------------------------------------------------------------------

# --- paragraph 34: replace text, add superscript "th" ---
$d = $word.ActiveDocument
$p34 = $d.Paragraphs.Item(34)
$full34 = $d.Range($p34.Range.Start, $p34.Range.End)
$full34.Text = "Faculties/CIO team will upload question paper by 24th November 2023"

$d = $word.ActiveDocument
$p34 = $d.Paragraphs.Item(34)
$text34 = $p34.Range.Text
$thIdx34 = $text34.IndexOf("24th") + 2
$thStart34 = $p34.Range.Start + $thIdx34
$thRange34 = $d.Range($thStart34, $thStart34 + 2)
$thRange34.Font.Superscript = $true

# --- paragraph 35: straightforward text replacement ---
$d = $word.ActiveDocument
$p35 = $d.Paragraphs.Item(35)
$full35 = $d.Range($p35.Range.Start, $p35.Range.End)
$full35.Text = "College Doors will create login credential for each faculties."

# --- paragraph 36: split into two list paragraphs ---
# First, insert a new (empty) list paragraph right after paragraph 36,
# inheriting the same ListParagraph / numbering formatting.
$d = $word.ActiveDocument
$p36 = $d.Paragraphs.Item(36)
$p36.Range.InsertParagraphAfter()

# Replace paragraph 36's text (this also removes its old "th" superscript
# run and trailing date text, since the whole run set is being replaced).
$d = $word.ActiveDocument
$p36 = $d.Paragraphs.Item(36)
$full36 = $d.Range($p36.Range.Start, $p36.Range.End)
$full36.Text = "We will give names of students and phone number to College Doors"

# Populate the newly-inserted (37th) paragraph.
$d = $word.ActiveDocument
$p37 = $d.Paragraphs.Item(37)
$p37.Range.Text = "College Doors will register students" + [char]0x2019 + " names in database by 24th November 2023"

# Force the run split between "...in database" and " by 24" (both pieces
# keep identical formatting, matching how the author typed it in two goes).
$d = $word.ActiveDocument
$p37 = $d.Paragraphs.Item(37)
$text37 = $p37.Range.Text
$splitIdx37 = $text37.IndexOf(" by 24")
$splitPos37 = $p37.Range.Start + $splitIdx37
$afterSplit37 = $d.Range($splitPos37, $p37.Range.End)
$afterSplit37.Bold = 1
$afterSplit37.Bold = 0

# Make the "th" in "24th" superscript.
$d = $word.ActiveDocument
$p37 = $d.Paragraphs.Item(37)
$text37 = $p37.Range.Text
$thIdx37 = $text37.IndexOf("24th") + 2
$thStart37 = $p37.Range.Start + $thIdx37
$thRange37 = $d.Range($thStart37, $thStart37 + 2)
$thRange37.Font.Superscript = $true
